$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.492.36"
$ws.Range("D2").Style = "Normal"

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.731.12"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.69%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.08%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.83"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.06%  "

# Row 6
$ws.Range("E6").Value = "  +0.03%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4883"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.22%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2670"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -1.00%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06218"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.76%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.733.02"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.61%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07067"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.86%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "15.67"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.30%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.645"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.79%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.6091"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.50%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "77.36"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.09%  "

# Row 16
$ws.Range("E16").Value = "  +0.03%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.484.21"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.43%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "1.001"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.07%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007186"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +4.10%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "11.50"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.50%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.956.33"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.61%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.523"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.00%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.771"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.25%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.252"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.18%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "139.23"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.03%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "15.43"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.36%  "

# Row 27
$ws.Range("E27").Value = "  -2.26%  "

# Row 28
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.404"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.82%  "

# Row 29
$ws.Range("B29").Value = "BitcoinCash"
$ws.Range("C29").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "107.93"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.99%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.970"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.16%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08050"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.02%  "

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.689"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.44%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04576"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -0.29%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.615"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.03%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.006"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.51%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6390"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.27%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9019"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -3.32%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.018"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +1.37%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.398"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -1.66%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.002"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.22%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.01508"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.27%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "101.40"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -10.82%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.443"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.33%  "

# Row 44
$ws.Range("E44").Value = "  -0.57%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "6.947"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.71%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1184"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.66%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05387"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.95%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "30.54"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.80%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.809"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -1.83%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.248"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.22%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3411"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.27%  "
